# This script refreshes the market-price / profit columns (H:N) of the
# per-job Leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to
# match the latest scheduled market-data pull. Only specific rows whose
# underlying item prices moved are touched; all other cells are untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2: Mercury Rising / Quicksilver
$ws.Range("H2").Value = 930
$ws.Range("I2").Value = 1000
$ws.Range("K2").Value = 1000
$ws.Range("M2").Value = -887

# Row 28: The Writing Is Not on the Wall / Enchanted Silver Ink
$ws.Range("H28").Value = 311.3846
$ws.Range("I28").Value = 192.44444
$ws.Range("K28").Value = 192.44444
$ws.Range("M28").Value = 292.55556

# Row 29: Dripping with Venom / Weak Blinding Potion
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

# Row 40: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 1013.3333
$ws.Range("I40").Value = 784
$ws.Range("K40").Value = 784
$ws.Range("M40").Value = -609

# Row 74: Adhesive of Antipathy / Wing Glue
$ws.Range("H74").Value = 20837944
$ws.Range("I74").Value = 4000
$ws.Range("J74").Value = 31254916
$ws.Range("K74").Value = 4000
$ws.Range("L74").Value = 31254916
$ws.Range("M74").Value = -3064
$ws.Range("N74").Value = -31256788

# Row 77: It's Gonna Grow Back (L) / Wing Glue
$ws.Range("H77").Value = 20837944
$ws.Range("I77").Value = 4000
$ws.Range("J77").Value = 31254916
$ws.Range("K77").Value = 20000
$ws.Range("L77").Value = 156274580
$ws.Range("M77").Value = -15320
$ws.Range("N77").Value = -156283940

# Row 129: Practical Command / Commanding Craftsman's Draught
$ws.Range("H129").Value = 304001.28
$ws.Range("J129").Value = 358237.22
$ws.Range("L129").Value = 1074711.66
$ws.Range("N129").Value = -1084711.66

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 21741348
$ws.Range("I132").Value = 23811822
$ws.Range("J132").Value = 1376.5
$ws.Range("K132").Value = 71435466
$ws.Range("L132").Value = 4129.5
$ws.Range("M132").Value = -71432936
$ws.Range("N132").Value = -9189.5

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 105197.16
$ws.Range("I137").Value = 140078.94
$ws.Range("J137").Value = 4040
$ws.Range("K137").Value = 420236.82
$ws.Range("L137").Value = 12120
$ws.Range("M137").Value = -417686.82
$ws.Range("N137").Value = -17220

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 15823.898
$ws.Range("I32").Value = 11197
$ws.Range("J32").Value = 24846.35
$ws.Range("K32").Value = 11197
$ws.Range("L32").Value = 24846.35
$ws.Range("M32").Value = -10910
$ws.Range("N32").Value = -25420.35

# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 4606.316
$ws.Range("I45").Value = 4599.385
$ws.Range("K45").Value = 4599.385
$ws.Range("M45").Value = -4222.385

# Row 60: Booty Call / Cobalt-plated Jackboots
$ws.Range("H60").Value = 18000
$ws.Range("I60").Value = 18000
$ws.Range("J60").Value = 18000
$ws.Range("K60").Value = 18000
$ws.Range("L60").Value = 18000
$ws.Range("M60").Value = -17267
$ws.Range("N60").Value = -19466

# Row 70: Pan That Laid the Golden Egg / Frypan Caliente
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("N70").Value = 0

# Row 73: Skillet with Fire (L) / Frypan Caliente
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("N73").Value = 0

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 30304864
$ws.Range("I74").Value = 50000800
$ws.Range("J74").Value = 3423.077
$ws.Range("K74").Value = 50000800
$ws.Range("L74").Value = 3423.077
$ws.Range("M74").Value = -49999926
$ws.Range("N74").Value = -5171.077

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 30304864
$ws.Range("I77").Value = 50000800
$ws.Range("J77").Value = 3423.077
$ws.Range("K77").Value = 250004000
$ws.Range("L77").Value = 17115.385
$ws.Range("M77").Value = -249999632
$ws.Range("N77").Value = -25851.385

# Row 128: Heading toward Bankruptcy / Manganese Helm of the Falling Dragon
$ws.Range("H128").Value = 33333.332
$ws.Range("J128").Value = 33333.332
$ws.Range("L128").Value = 33333.332
$ws.Range("N128").Value = -43293.332

$ws = $wb.Worksheets.Item("BSM")
# Row 25: Tools of the Trade / Iron Doming Hammer
$ws.Range("H25").Value = 491.33334
$ws.Range("I25").Value = 491.33334
$ws.Range("K25").Value = 491.33334
$ws.Range("M25").Value = -256.33334

# Row 82: Spirituality Inspector / Titanium Lump Hammer
$ws.Range("H82").Value = 27511.125
$ws.Range("I82").Value = 7816
$ws.Range("K82").Value = 7816
$ws.Range("M82").Value = -7433

# Row 85: The Clamor for Hammers (L) / Titanium Lump Hammer
$ws.Range("H85").Value = 27511.125
$ws.Range("I85").Value = 7816
$ws.Range("K85").Value = 7816
$ws.Range("M85").Value = -6490

# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 1969.3334
$ws.Range("I86").Value = 1813.1052
$ws.Range("J86").Value = 3453.5
$ws.Range("K86").Value = 1813.1052
$ws.Range("L86").Value = 3453.5
$ws.Range("M86").Value = -690.1052
$ws.Range("N86").Value = -5699.5

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 1969.3334
$ws.Range("I89").Value = 1813.1052
$ws.Range("J89").Value = 3453.5
$ws.Range("K89").Value = 9065.526
$ws.Range("L89").Value = 17267.5
$ws.Range("M89").Value = -3449.526
$ws.Range("N89").Value = -28499.5

# Row 94: High Steal / High Steel Nugget
$ws.Range("H94").Value = 830.85297
$ws.Range("I94").Value = 340.82608
$ws.Range("K94").Value = 340.82608
$ws.Range("M94").Value = 110.17392

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 4094.4167
$ws.Range("I107").Value = 2998
$ws.Range("K107").Value = 2998
$ws.Range("M107").Value = -1078

$ws = $wb.Worksheets.Item("CRP")
# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 5589.5557
$ws.Range("J62").Value = 7668.6665
$ws.Range("L62").Value = 7668.6665
$ws.Range("N62").Value = -8916.666499999999

# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 5589.5557
$ws.Range("J65").Value = 7668.6665
$ws.Range("L65").Value = 38343.3325
$ws.Range("N65").Value = -44583.3325

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 4047.9167
$ws.Range("I99").Value = 3105.8823
$ws.Range("J99").Value = 6335.7144
$ws.Range("K99").Value = 3105.8823
$ws.Range("L99").Value = 6335.7144
$ws.Range("M99").Value = -1607.8823
$ws.Range("N99").Value = -9331.714400000001

# Row 105: Zelkova, My Love / Zelkova Lumber
$ws.Range("H105").Value = 6511
$ws.Range("J105").Value = 6511
$ws.Range("L105").Value = 6511
$ws.Range("N105").Value = -10005

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 4047.9167
$ws.Range("I126").Value = 3105.8823
$ws.Range("J126").Value = 6335.7144
$ws.Range("K126").Value = 9317.6469
$ws.Range("L126").Value = 19007.1432
$ws.Range("M126").Value = -6847.6469
$ws.Range("N126").Value = -23947.1432

$ws = $wb.Worksheets.Item("CUL")
# Row 23: Sweet Smell of Success / Lavender Oil
$ws.Range("H23").Value = 261.32
$ws.Range("I23").Value = 16.666666
$ws.Range("J23").Value = 294.68182
$ws.Range("K23").Value = 49.999998
$ws.Range("L23").Value = 884.04546
$ws.Range("M23").Value = 185.000002
$ws.Range("N23").Value = -1354.04546

# Row 113: Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 851.25
$ws.Range("I113").Value = 604.1111
$ws.Range("J113").Value = 1053.4546
$ws.Range("K113").Value = 1812.3333
$ws.Range("L113").Value = 3160.3638
$ws.Range("M113").Value = 357.6667000000002
$ws.Range("N113").Value = -7500.3638

# Row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 1445.1923
$ws.Range("I122").Value = 489
$ws.Range("J122").Value = 1524.875
$ws.Range("K122").Value = 4401
$ws.Range("L122").Value = 13723.875
$ws.Range("M122").Value = -1951
$ws.Range("N122").Value = -18623.875

# Row 123: Topping Up the Pot / Zurek
$ws.Range("H123").Value = 2888.3333
$ws.Range("J123").Value = 5595
$ws.Range("L123").Value = 16785
$ws.Range("N123").Value = -21685

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 688.28864
$ws.Range("I131").Value = 401.66666
$ws.Range("J131").Value = 728.7529
$ws.Range("K131").Value = 1204.99998
$ws.Range("L131").Value = 2186.2587
$ws.Range("M131").Value = 3835.00002
$ws.Range("N131").Value = -12266.2587

# Row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 651.2857
$ws.Range("I132").Value = 693.1667
$ws.Range("J132").Value = 400
$ws.Range("K132").Value = 6238.5003
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -3708.5003
$ws.Range("N132").Value = -8660

$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Range("H97").Value = 3018.5
$ws.Range("I97").Value = 3400
$ws.Range("J97").Value = 2827.75
$ws.Range("K97").Value = 3400
$ws.Range("L97").Value = 2827.75
$ws.Range("M97").Value = -2904
$ws.Range("N97").Value = -3819.75

# Row 130: Planisphere to Paper / Chondrite Magitek Planisphere
$ws.Range("H130").Value = 52109.816
$ws.Range("J130").Value = 52109.816
$ws.Range("L130").Value = 52109.816
$ws.Range("N130").Value = -62149.816

# Row 131: Star Athletes / Star Quartz Wristband of Aiming
$ws.Range("H131").Value = 49664.285
$ws.Range("J131").Value = 49664.285
$ws.Range("L131").Value = 49664.285
$ws.Range("N131").Value = -59744.285

$ws = $wb.Worksheets.Item("LTW")
# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 1088.1177
$ws.Range("I46").Value = 990.36365
$ws.Range("J46").Value = 1267.3334
$ws.Range("K46").Value = 990.36365
$ws.Range("L46").Value = 1267.3334
$ws.Range("M46").Value = -802.36365
$ws.Range("N46").Value = -1643.3334

# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Range("H55").Value = 268.38095
$ws.Range("I55").Value = 261.125
$ws.Range("K55").Value = 261.125
$ws.Range("M55").Value = -88.125

# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 2567
$ws.Range("J68").Value = 3167.3333
$ws.Range("L68").Value = 3167.3333
$ws.Range("N68").Value = -4665.3333

# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 2567
$ws.Range("J71").Value = 3167.3333
$ws.Range("L71").Value = 15836.6665
$ws.Range("N71").Value = -23324.6665

# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 4221.2
$ws.Range("I82").Value = 2366.6667
$ws.Range("K82").Value = 2366.6667
$ws.Range("M82").Value = -2005.6667

# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 4221.2
$ws.Range("I85").Value = 2366.6667
$ws.Range("K85").Value = 2366.6667
$ws.Range("M85").Value = -1118.6667

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 2852238.5
$ws.Range("I122").Value = 4987367.5
$ws.Range("K122").Value = 14962102.5
$ws.Range("M122").Value = -14959652.5

$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke / Rainbow Cloth
$ws.Range("H62").Value = 3994.4
$ws.Range("I62").Value = 3994.4
$ws.Range("K62").Value = 3994.4
$ws.Range("M62").Value = -3370.4

# Row 65: Desperate for Diversionaries (L) / Rainbow Cloth
$ws.Range("H65").Value = 3994.4
$ws.Range("I65").Value = 3994.4
$ws.Range("K65").Value = 19972
$ws.Range("M65").Value = -16852

# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 1929
$ws.Range("I122").Value = 1833.8334
$ws.Range("K122").Value = 5501.5002
$ws.Range("M122").Value = -3051.5002

# Row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 2714.2144
$ws.Range("I126").Value = 2277.6667
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 6833.000100000001
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -4363.000100000001
$ws.Range("N126").Value = -15440
